# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 125 for
# "Feria Lagunitas de Puerto Montt - Pomelo", pushing the existing
# rows 125:147 down to 126:148.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 125 (shifts rows 125-147 -> 126-148)
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row with this week's data
$ws.Range("A125").Value = 4
$ws.Range("B125").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C125").Value = "Los Lagos"
$ws.Range("D125").Value = 44474
$ws.Range("E125").Value = 10
$ws.Range("F125").Value = "Fruta"
$ws.Range("G125").Value = 100102
$ws.Range("H125").Value = "Cítricos"
$ws.Range("I125").Value = 100102006
$ws.Range("J125").Value = "Pomelo"
$ws.Range("K125").Value = "Start Ruby"
$ws.Range("L125").Value = "Primera"
$ws.Range("M125").Value = 100
$ws.Range("N125").Value = 12000
$ws.Range("O125").Value = 12000
$ws.Range("P125").Value = 12000
$ws.Range("Q125").Value = "$/caja 14 kilos empedrada"
$ws.Range("R125").Value = "Región de O'Higgins"
$ws.Range("S125").Value = 857
$ws.Range("T125").Value = 14

# Make sure the D125 cell keeps the same date-formatted style ("s=2")
# that every other date cell in column D uses.
$ws.Range("D125").NumberFormat = $ws.Range("D126").NumberFormat
